$d = $word.ActiveDocument

# 1. Update the creation-date timestamp everywhere it appears (title block + header).
$d.Content.Find.Execute("2021-10-24 17:28", $true, $false, $false, $false, $false, $true, 1, $false, "2021-10-24 19:10", 2) | Out-Null

$sec = $d.Sections.First
$hdr = $sec.Headers(1)
$hdr.Range.Find.Execute("2021-10-24 17:28", $true, $false, $false, $false, $false, $true, 1, $false, "2021-10-24 19:10", 2) | Out-Null

# 2. Append the new "Opportunity" section (heading, intro paragraph, numbered list items)
#    right after the Introduction paragraph / before the section break.
function Add-OppPara([string]$style, [string]$text) {
    $p = $d.Paragraphs.Last
    $r = $p.Range
    $r.Collapse(0)
    $r.InsertParagraphAfter() | Out-Null
    $np = $d.Paragraphs.Last
    $np.Style = $style
    $np.Range.Text = $text
}

$lineBreak = [char]11

Add-OppPara "Heading1" ("Opportunity")
Add-OppPara "Normal" ("Overall the two phases of the study were well paired. Essentially both the first and second phases illustrated that product managers need tooling that enable them to discern the whys behind an effort leadning to a distillation of whats needed to build an offer.  Thus there is an opportunity to build a product enabling product managers (and potentially customer success managers) to reveal insights from interactions with users, partners and competitors powering efforts ranging from from product modernization to new product introduction. Further, competitive insights surfaced in the second phase related to Aha!.  These insights showed that Aha! has yet to tackle features that automatically and systematically use interactions to reveal the whys and whats behind product roadmaps.  (Note that Aha! is highly relevant to sense competitive intelligence because they are the leader in market category of Product Management and Roadmapping tooling.)  While these two phases paired well, a weakness was obvious when the two phases were combined:  An intersection between them did not clearly surface.  Therefore a third phase was performed to validate and look for clear couplings between the first two phases.  What follows are some of the key opportunities discovered through these three phases followed by discrete sections relating systematically uncovered key themes, snippets associated to key themes, and finally abstracts for related interactions.")
Add-OppPara "ListNumber" ("Formalize Product Relationship Management - Tooling and associated process is needed to enable the build out of the community around the product management team.")
Add-OppPara "ListNumber" ("Outlive the Product Manager across the entire lifecycle - Product managers are not always present throughout an entire program lifecycle or may leave the company making it essential that source materials and decisioning reasoning stands alone.")
Add-OppPara "ListNumber" ("Integrate critical stakeholders via tooling - Product Managers are a key part of an overall program, but they do no live on an island making tooling access for a diverse set of stakeholders required.")
Add-OppPara "ListNumber" ("Reduce the time and effort of product research and feedback correlation - Discerning the core whys and whats of any offer is super critical, but the path to get there is often slow and intransparent.  Therefore, tooling should" + $lineBreak + "drive speed, improve transparency and reduce work burden.")
Add-OppPara "ListNumber" ("Forward and backward Traceability from problem to solution - While modern product management and roadmapping tooling facilitates process transparency, getting to clear and key whys and whats is frequently opaque and" + $lineBreak + "untraceable.  Clearly, revealing the path from problem identication, the whys, to problem resolution, the whats, is a key opportunity for tooling.")
Add-OppPara "ListNumber" ("Intelligent Information Integration by connecting Productivity, CRM, Support, PM tools - New tooling cannot exist in an island therefore any opportunity requires integration into a user/customer ecosystem.")
Add-OppPara "ListNumber" ("Visibility and Reporting for relevant stakeholders by themes, products and customers - Beyond kicking off work with engineering many stakeholders want to understand how whys, encoded in key themes, are being" + $lineBreak + "resolved.  For example customer success managers will need to know how their customers have affected the roadmap, customers themselves would like to understand their level of influence, and" + $lineBreak + "marketing teams will want to map features to key user pain points.  This means an offer should enable all interested stakeholders to ask and answer key questions beyond the what is needed to drive a roadmap.")
Add-OppPara "ListNumber" ("PRFAQ - In the third round of interviews it has become obvious that the output format should look more like an Amazon PR-FAQ which stands for Press Release and" + $lineBreak + "Frequently Asked Questions.  Therefore as we progress towards MVP it will be required to change the report format to PR-FAQ.")
